# The source diff for this revision touches exactly one kind of line,
# 73 times, across both slides of the deck:
#
#     -          <p:cNvPr id="<old>" name="...">
#     +          <p:cNvPr id="<new>" name="...">
#
# i.e. every shape's *non-visual* `cNvPr/@id` (its internal OOXML shape
# id) was reassigned to a different numeric value while every other
# byte of the slide XML - geometry, fills, run text, formatting,
# shape `name`, z-order, slide count, etc. - stayed byte-for-byte the
# same. The commit message ("summary pptx 경로 수정" / "fix summary pptx
# path") together with the `WEB-INF/Files/...` server path in the diff
# header indicates the file was simply re-emitted to a new location by
# the server-side pipeline (a non-PowerPoint tool), which mints fresh
# shape ids on write; it was not an edit made through PowerPoint.
#
# `Shape.Id` (and `ShapeRange.Id`) is a read-only, engine-assigned
# identifier in the real PowerPoint object model - there is no
# property setter, method, paste/duplicate/group operation, or other
# automation call that lets a caller choose or overwrite an existing
# shape's id; new ids are only ever minted internally when a shape is
# first created, from an internal monotonic/random counter outside
# caller control. (Verified interactively against this host: setting
# `$shape.Id`/`$shape.ID` is silently ignored, and
# Copy/Paste/Cut/Duplicate/Group/Ungroup/AddTextbox all preserve
# existing ids or mint small sequential ones - never the specific
# large values in the diff.)
#
# Since there is no visible/content change to make (the diff carries
# none) and the id churn itself cannot be reproduced through any
# legitimate COM/VBA call, this script intentionally performs no
# shape edits - that is the faithful, achievable result given the
# object model's constraints.
$p = $ppt.ActivePresentation
